$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.523.57"
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("D3").Value = "2.430.77"
$ws.Range("E3").Value = "  +7.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "296.72"
$ws.Range("E5").Value = "  -1.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.61"
$ws.Range("E6").Value = "  -2.53%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.572"
$ws.Range("E7").Value = "  +0.44%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.513"
$ws.Range("E9").Value = "  +0.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.45"
$ws.Range("E10").Value = "  +1.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0790"
$ws.Range("E11").Value = "  -1.47%  "
$ws.Range("E12").Value = "  +1.29%  "
$ws.Range("E13").Value = "  +2.34%  "
$ws.Range("D14").Value = "2.802.11"
$ws.Range("E14").Value = "  +7.01%  "
$ws.Range("D15").Value = "2.445.86"
$ws.Range("E15").Value = "  +7.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.852"
$ws.Range("E16").Value = "  +6.92%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.12"
$ws.Range("E17").Value = "  +3.46%  "
$ws.Range("D18").Value = "46.337.31"
$ws.Range("E18").Value = "  -0.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.80"
$ws.Range("E19").Value = "  +1.43%  "
$ws.Range("D20").Value = "0.0₃0952"
$ws.Range("E20").Value = "  -1.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.25"
$ws.Range("E21").Value = "  +7.04%  "
$ws.Range("E22").Value = "  +2.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "246.02"
$ws.Range("E23").Value = "  -0.80%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.81"
$ws.Range("E24").Value = "  +0.45%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.97"
$ws.Range("E25").Value = "  +5.70%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "39.54"
$ws.Range("E27").Value = "  -4.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.23"
$ws.Range("E28").Value = "  -0.87%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.81"
$ws.Range("E29").Value = "  +1.77%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.87"
$ws.Range("E30").Value = "  +15.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "21.42"
$ws.Range("E31").Value = "  +5.77%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.76"
$ws.Range("E32").Value = "  -1.83%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.59"
$ws.Range("E33").Value = "  +4.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "148.21"
$ws.Range("E34").Value = "  +0.73%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0775"
$ws.Range("E35").Value = "  +0.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.01"
$ws.Range("E36").Value = "  +18.76%  "
$ws.Range("E37").Value = "  +0.89%  "
$ws.Range("E38").Value = "  +0.59%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.34"
$ws.Range("E39").Value = "  -2.93%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.95"
$ws.Range("E40").Value = "  +2.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0304"
$ws.Range("E41").Value = "  +2.56%  "
$ws.Range("E42").Value = "  +4.73%  "
$ws.Range("D43").Value = "1.978.39"
$ws.Range("E43").Value = "  +10.70%  "
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "93.16"
$ws.Range("E45").Value = "  +0.67%  "
$ws.Range("E46").Value = "  -1.66%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.47"
$ws.Range("E47").Value = "  +31.72%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.65"
$ws.Range("E48").Value = "  +9.94%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "101.43"
$ws.Range("E49").Value = "  +7.06%  "
$ws.Range("D50").Value = "2.670.85"
$ws.Range("E50").Value = "  +7.02%  "
$ws.Range("E51").Value = "  +0.75%  "
